$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 286; this shifts the existing rows 286..311
# down to 287..312 (values, styles and the sheet's dimension are all
# updated automatically by Excel's row-insert semantics).
$ws.Rows("286").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A286").Value = 10
$ws.Range("B286").Value = "Vega Modelo de Temuco"
$ws.Range("C286").Value = "La Araucanía"
$ws.Range("D286").Value = 44461
$ws.Range("E286").Value = 9
$ws.Range("F286").Value = 100112028
$ws.Range("G286").Value = "Sandia"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 280
$ws.Range("K286").Value = 1100
$ws.Range("L286").Value = 1100
$ws.Range("M286").Value = 1100
$ws.Range("N286").Value = "$/kilo (volumen en unidades)"
$ws.Range("O286").Value = "Perú"
$ws.Range("P286").Value = 1100
$ws.Range("Q286").Value = 1
$ws.Range("R286").Value = "Hortaliza"
